$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1900
$ws.Range("I29").Value = 233.33333
$ws.Range("J29").Value = 2316.6667
$ws.Range("K29").Value = 699.99999
$ws.Range("L29").Value = 6950.000100000001
$ws.Range("M29").Value = -418.99999
$ws.Range("N29").Value = -7512.000100000001
$ws.Range("H40").Value = 1971.5714
$ws.Range("J40").Value = 2080.4
$ws.Range("L40").Value = 2080.4
$ws.Range("N40").Value = -2430.4
$ws.Range("H62").Value = 22225604
$ws.Range("I62").Value = 27781528
$ws.Range("J62").Value = 1906
$ws.Range("K62").Value = 27781528
$ws.Range("L62").Value = 1906
$ws.Range("M62").Value = -27780904
$ws.Range("N62").Value = -3154
$ws.Range("H65").Value = 22225604
$ws.Range("I65").Value = 27781528
$ws.Range("J65").Value = 1906
$ws.Range("K65").Value = 138907640
$ws.Range("L65").Value = 9530
$ws.Range("M65").Value = -138904520
$ws.Range("N65").Value = -15770
$ws.Range("H82").Value = 1725
$ws.Range("I82").Value = 800
$ws.Range("K82").Value = 2400
$ws.Range("M82").Value = -1994
$ws.Range("H85").Value = 1725
$ws.Range("I85").Value = 800
$ws.Range("K85").Value = 2400
$ws.Range("M85").Value = -996
$ws.Range("H98").Value = 1810.3833
$ws.Range("I98").Value = 1900.6727
$ws.Range("J98").Value = 817.2
$ws.Range("K98").Value = 1900.6727
$ws.Range("L98").Value = 817.2
$ws.Range("M98").Value = -402.6727000000001
$ws.Range("N98").Value = -3813.2
$ws.Range("H122").Value = 1810.3833
$ws.Range("I122").Value = 1900.6727
$ws.Range("J122").Value = 817.2
$ws.Range("K122").Value = 5702.0181
$ws.Range("L122").Value = 2451.6
$ws.Range("M122").Value = -3252.0181
$ws.Range("N122").Value = -7351.6
$ws.Range("H137").Value = 1503.6578
$ws.Range("I137").Value = 1001.6
$ws.Range("J137").Value = 2061.5
$ws.Range("K137").Value = 3004.8
$ws.Range("L137").Value = 6184.5
$ws.Range("M137").Value = -454.8000000000002
$ws.Range("N137").Value = -11284.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = $null
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = $null
$ws.Range("H102").Value = 16668776
$ws.Range("I102").Value = 20835226
$ws.Range("K102").Value = 20835226
$ws.Range("M102").Value = -20833604

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 2950
$ws.Range("J17").Value = 900
$ws.Range("L17").Value = 900
$ws.Range("N17").Value = -1244
$ws.Range("H24").Value = 1800
$ws.Range("J24").Value = 1800
$ws.Range("L24").Value = 1800
$ws.Range("N24").Value = -2270
$ws.Range("H25").Value = 5016
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 5016
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 5016
$ws.Range("M25").Value = $null
$ws.Range("N25").Value = -5486
$ws.Range("H86").Value = 2944.8333
$ws.Range("I86").Value = 3131.6875
$ws.Range("J86").Value = 1450
$ws.Range("K86").Value = 3131.6875
$ws.Range("L86").Value = 1450
$ws.Range("M86").Value = -2008.6875
$ws.Range("N86").Value = -3696
$ws.Range("H89").Value = 2944.8333
$ws.Range("I89").Value = 3131.6875
$ws.Range("J89").Value = 1450
$ws.Range("K89").Value = 15658.4375
$ws.Range("L89").Value = 7250
$ws.Range("M89").Value = -10042.4375
$ws.Range("N89").Value = -18482
$ws.Range("H134").Value = 2223.4092
$ws.Range("I134").Value = 1643.9445
$ws.Range("J134").Value = 4831
$ws.Range("K134").Value = 4931.833500000001
$ws.Range("L134").Value = 14493
$ws.Range("M134").Value = -2396.833500000001
$ws.Range("N134").Value = -19563

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 247.2
$ws.Range("I7").Value = 152
$ws.Range("J7").Value = 437.6
$ws.Range("K7").Value = 152
$ws.Range("L7").Value = 437.6
$ws.Range("M7").Value = -39
$ws.Range("N7").Value = -663.6
$ws.Range("H19").Value = 711.2
$ws.Range("I19").Value = 137.5
$ws.Range("J19").Value = 3006
$ws.Range("K19").Value = 137.5
$ws.Range("L19").Value = 3006
$ws.Range("M19").Value = 32.5
$ws.Range("N19").Value = -3346
$ws.Range("H23").Value = 10000
$ws.Range("J23").Value = 10000
$ws.Range("L23").Value = 10000
$ws.Range("N23").Value = -10480
$ws.Range("H24").Value = 711.2
$ws.Range("I24").Value = 137.5
$ws.Range("J24").Value = 3006
$ws.Range("K24").Value = 137.5
$ws.Range("L24").Value = 3006
$ws.Range("M24").Value = 32.5
$ws.Range("N24").Value = -3346
$ws.Range("H27").Value = 10000
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10384
$ws.Range("H31").Value = 1546.8966
$ws.Range("I31").Value = 1146.5652
$ws.Range("J31").Value = 3081.5
$ws.Range("K31").Value = 1146.5652
$ws.Range("L31").Value = 3081.5
$ws.Range("M31").Value = -851.5652
$ws.Range("N31").Value = -3671.5
$ws.Range("H34").Value = 1546.8966
$ws.Range("I34").Value = 1146.5652
$ws.Range("J34").Value = 3081.5
$ws.Range("K34").Value = 1146.5652
$ws.Range("L34").Value = 3081.5
$ws.Range("M34").Value = -944.5652
$ws.Range("N34").Value = -3485.5
$ws.Range("H141").Value = 768264.75
$ws.Range("J141").Value = 768264.75
$ws.Range("L141").Value = 768264.75
$ws.Range("N141").Value = -778624.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1113.6364
$ws.Range("I17").Value = 1113.6364
$ws.Range("K17").Value = 3340.9092
$ws.Range("M17").Value = -3171.9092
$ws.Range("H58").Value = 3320
$ws.Range("J58").Value = 3320
$ws.Range("L58").Value = 9960
$ws.Range("N58").Value = -10216
$ws.Range("H131").Value = 26319908
$ws.Range("J131").Value = 4837.4375
$ws.Range("L131").Value = 14512.3125
$ws.Range("N131").Value = -24592.3125
$ws.Range("H138").Value = 1945.08
$ws.Range("I138").Value = 1676.75
$ws.Range("J138").Value = 2422.111
$ws.Range("K138").Value = 5030.25
$ws.Range("L138").Value = 7266.333
$ws.Range("M138").Value = 109.75
$ws.Range("N138").Value = -17546.333
$ws.Range("H141").Value = 2601.0833
$ws.Range("I141").Value = 1661.8182
$ws.Range("K141").Value = 4985.4546
$ws.Range("M141").Value = 194.5454

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 6900
$ws.Range("J29").Value = 6900
$ws.Range("L29").Value = 6900
$ws.Range("N29").Value = -7480
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = $null
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = $null
$ws.Range("H102").Value = 1292.5454
$ws.Range("I102").Value = 1071.8
$ws.Range("K102").Value = 1071.8
$ws.Range("M102").Value = 550.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 800
$ws.Range("I4").Value = 800
$ws.Range("K4").Value = 800
$ws.Range("M4").Value = -687
$ws.Range("H28").Value = 800
$ws.Range("I28").Value = 800
$ws.Range("K28").Value = 800
$ws.Range("M28").Value = -568
$ws.Range("H37").Value = 800
$ws.Range("I37").Value = 800
$ws.Range("K37").Value = 800
$ws.Range("M37").Value = -693
$ws.Range("H68").Value = 1988.2916
$ws.Range("I68").Value = 1929.5238
$ws.Range("J68").Value = 2399.6667
$ws.Range("K68").Value = 1929.5238
$ws.Range("L68").Value = 2399.6667
$ws.Range("M68").Value = -1180.5238
$ws.Range("N68").Value = -3897.6667
$ws.Range("H71").Value = 1988.2916
$ws.Range("I71").Value = 1929.5238
$ws.Range("J71").Value = 2399.6667
$ws.Range("K71").Value = 9647.618999999999
$ws.Range("L71").Value = 11998.3335
$ws.Range("M71").Value = -5903.618999999999
$ws.Range("N71").Value = -19486.3335
$ws.Range("H135").Value = 36555.57
$ws.Range("J135").Value = 36555.57
$ws.Range("L135").Value = 36555.57
$ws.Range("N135").Value = -46695.57
$ws.Range("H136").Value = 11603.7
$ws.Range("I136").Value = 15362.429
$ws.Range("J136").Value = 2833.3333
$ws.Range("K136").Value = 46087.287
$ws.Range("L136").Value = 8499.999899999999
$ws.Range("M136").Value = -43537.287
$ws.Range("N136").Value = -13599.9999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4256.0625
$ws.Range("I132").Value = 6274.75
$ws.Range("K132").Value = 18824.25
$ws.Range("M132").Value = -16294.25
$ws.Range("H136").Value = 990.1875
$ws.Range("I136").Value = 967.0476
$ws.Range("J136").Value = 1034.3636
$ws.Range("K136").Value = 2901.1428
$ws.Range("L136").Value = 3103.0908
$ws.Range("M136").Value = -351.1428000000001
$ws.Range("N136").Value = -8203.0908
